$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "324.84" or
# "1.929.72" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.915.25"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.907.94"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "324.84"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").Value = "0.3811"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").Value = "0.07717"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "0.9796"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "22.05"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").Value = "1.929.72"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "5.673"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "6.938"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").Value = "0.07056"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "83.73"
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("D18").Value = "0.000009458"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D19").Value = "16.63"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "28.907.89"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").Value = "5.320"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D23").Value = "10.92"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "158.66"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").Value = "19.04"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "5.665"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").Value = "1.869"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").Value = "0.09290"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").Value = "0.8626"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "5.079"
$ws.Range("D33").Value = "1.248"
$ws.Range("E33").Value = "  -4.40%  "
$ws.Range("D34").Value = "3.076"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "0.05713"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("D36").Value = "1.158"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").Value = "7.420"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("D40").Value = "0.5486"
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("D42").Value = "2.884"
$ws.Range("E42").Value = "  +7.08%  "
$ws.Range("D43").Value = "9.310"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("D45").Value = "2.126"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("D46").Value = "11.19"
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("D47").Value = "0.06888"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").Value = "110.23"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "0.000002562"
$ws.Range("E50").Value = "  -16.07%  "
$ws.Range("D51").Value = "0.2858"
$ws.Range("E51").Value = "  -4.73%  "
